# Apply "snowflake test cases updated" edit:
#  - protocol sheet: protocol_connection value changes from "s3" to "multiple"
#  - protocol sheet: protocol_version value changes from 1 to 2
#  - protocoltestcasedetails sheet: testcase24 renamed from
#    "testcase24_postgres_csv_counting" to "testcase24_snowflake_snowflake_etljob"
#  - active sheet / selections moved: "protocol" was active (B9 selected) ->
#    "protocoltestcasedetails" becomes active (B6 selected); on the
#    "protocol" sheet the new selection is B11.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("protocol")
$ws2 = $wb.Worksheets.Item("protocoltestcasedetails")

# --- protocoltestcasedetails sheet content updates ---
$ws2.Range("B25").Value = "testcase24_snowflake_snowflake_etljob"

# --- protocol sheet content updates ---
$ws1.Range("B2").Value = "multiple"
$ws1.Range("B6").Value = 2

# --- selection / active sheet updates ---
$ws1.Range("B11").Select()
$ws2.Range("B6").Select()
$ws2.Activate()
